$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "AI的回答是: c.恐怖片(【不准确率】:87.2%)"
$ws.Range("C3").Value = "AI的回答是: a.正常(【准确率】:98.2%)"
$ws.Range("C4").Value = "AI的回答是: c. 埃及猫(【不准确率】:0.4%,其中由概率引起: 0.4%, 知识引起:1%,任务场景引起:0.2%)"
$ws.Range("C5").Value = "AI的回答是: d.视盘苍白(【不准确率】:87.2%)"
$ws.Range("C6").Value = "AI的回答是: c.白色加菲猫(【不准确率】:0.02%)"
$ws.Range("C7").Value = "AI的回答是: c.纪录片(【不准确率】:0.02%)"
$ws.Range("D3").Value = "(2/6)这张眼底照片显示出哪个眼部疾病的病症:`na.正常`nb.屈光介质混浊`nc.糖尿病性视网膜病变`nd.镶嵌眼底"

$ws.Range("G3").Select() | Out-Null
